$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45179 = 2023-09-10)
# that was bumped by one day (45180 = 2023-09-11) for every data row
# (rows 2 through 525).
$ws.Range("C2:C525").Value = 45180
